$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 400, shifting existing rows 400-425 down to 401-426.
$ws.Rows.Item(400).Insert()

# Populate the newly inserted row 400 with the new weekly record.
$ws.Cells.Item(400, 1).Value = 6
$ws.Cells.Item(400, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(400, 3).Value = "Metropolitana"
$ws.Cells.Item(400, 4).Value = 44516
$ws.Cells.Item(400, 4).Style = $ws.Cells.Item(401, 4).Style
$ws.Cells.Item(400, 4).NumberFormat = $ws.Cells.Item(401, 4).NumberFormat
$ws.Cells.Item(400, 5).Value = 13
$ws.Cells.Item(400, 6).Value = 100112003
$ws.Cells.Item(400, 7).Value = "Ajo"
$ws.Cells.Item(400, 8).Value = "Chino"
$ws.Cells.Item(400, 9).Value = "Primera"
$ws.Cells.Item(400, 10).Value = 2600
$ws.Cells.Item(400, 11).Value = 16500
$ws.Cells.Item(400, 12).Value = 17000
$ws.Cells.Item(400, 13).Value = 16712
$ws.Cells.Item(400, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(400, 15).Value = "China"
$ws.Cells.Item(400, 16).Value = 1671
$ws.Cells.Item(400, 17).Value = 10
$ws.Cells.Item(400, 18).Value = "Hortaliza"
